$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# positions sheet: MSD Capital, L.P. role now has an end date (they left the
# role) -- F2 = 12/1/2019 (serial 43800).
# ---------------------------------------------------------------------------
$positions = $wb.Worksheets.Item("positions")
$positions.Range("F2").Value = 43800

# ---------------------------------------------------------------------------
# projects sheet: restructure the header/data columns.
#   - column B used to hold the free-text project title, and column C held a
#     formula pulling the institution name from the positions sheet.
#   - now column B holds that institution-name formula, column C holds the
#     project "name" (the old title text) and a brand new column D holds an
#     "overview" tag for each project. Column E ("detail_1" per the header
#     row) keeps its existing meaning -- only one text value is corrected
#     (row 4) and two brand new rows of tag/name data are appended (row 5
#     gains an in_resume flag plus overview/detail tags).
# ---------------------------------------------------------------------------
$projects = $wb.Worksheets.Item("projects")

# Header row
$projects.Range("B1").Value = "institution"
$projects.Range("C1").Value = "name"
$projects.Range("D1").Value = "overview"

# Row 2 - Analytics Platform (MSD Capital, L.P.)
$projects.Range("C2").Value = "Analytics Platform"
$projects.Range("B2").Formula = "=positions!D$2"
$projects.Range("D2").Value = "Msd 1"

# Row 3 - Option Pricing Framework (MSD Capital, L.P.)
$projects.Range("C3").Value = "Option Pricing Framework"
$projects.Range("B3").Formula = "=positions!D$2"
$projects.Range("D3").Value = "trading"

# Row 4 - Independent Foreclosure Review (Promontory Financial Group)
$projects.Range("C4").Value = "Independent Foreclosure Review"
$projects.Range("B4").Formula = "=positions!D$3"
$projects.Range("D4").Value = "fa"
$projects.Range("E4").Value = "mortgage review"

# Row 5 - RightContent (RR Donnelley) - newly flagged as in_resume
$projects.Range("A4").Copy($projects.Range("A5"))
$projects.Range("C5").Value = "RightContent"
$projects.Range("B5").Formula = "=positions!D$4"
$projects.Range("D5").Value = "fsad"
$projects.Range("E5").Value = "Web App"

# ---------------------------------------------------------------------------
# Active tab: the "positions" sheet is now the one shown when the workbook
# is opened (was "projects").
# ---------------------------------------------------------------------------
$positions.Activate()
